$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D) updates: force text format so numeric-looking strings
# (e.g. trailing zeros like "3.500") are preserved exactly as text, matching
# the workbook convention where every Price cell is stored as inline text. ---
$priceCells = @("D2","D4","D5","D6","D8","D9","D10","D12","D13","D14","D15","D16","D17","D18","D20","D22","D23","D25","D26","D40","D41","D42","D43","D44","D45","D47","D48")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Assign new price values (as text)
$ws.Range("D2").Value = "245.46"
$ws.Range("D4").Value = "5.121"
$ws.Range("D5").Value = "0.05572"
$ws.Range("D6").Value = "6.498"
$ws.Range("D8").Value = "0.8171"
$ws.Range("D9").Value = "0.8443"
$ws.Range("D10").Value = "0.1344"
$ws.Range("D12").Value = "0.02889"
$ws.Range("D13").Value = "0.09376"
$ws.Range("D14").Value = "0.001529"
$ws.Range("D15").Value = "0.0005993"
$ws.Range("D16").Value = "0.006139"
$ws.Range("D17").Value = "3.500"
$ws.Range("D18").Value = "2.063"
$ws.Range("D20").Value = "0.03176"
$ws.Range("D22").Value = "3.744"
$ws.Range("D23").Value = "0.04714"
$ws.Range("D25").Value = "0.001251"
$ws.Range("D26").Value = "0.004639"
$ws.Range("D40").Value = "0.03661"
$ws.Range("D41").Value = "0.1366"
$ws.Range("D42").Value = "0.002631"
$ws.Range("D43").Value = "0.003381"
$ws.Range("D44").Value = "0.008316"
$ws.Range("D45").Value = "0.00005300"
$ws.Range("D47").Value = "0.1501"
$ws.Range("D48").Value = "0.002117"

# --- Coin name / link / volume-label updates (plain text, no numeric coercion risk) ---
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E15").Value = "14OneONE"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("E16").Value = "15TigerCashTCH"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("E17").Value = "16LEOLEO"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("E18").Value = "17BTSETokenBTSE"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E43").Value = "42KickTokenKICK"
